# Sync up for DMM capture
$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# The workbook starts with two sheets: "3_3" (active) and "1_8".
# Keep "3_3" (renaming it to "test"), drop "1_8" entirely.
$ws = $wb.Worksheets.Item("3_3")
$old = $wb.Worksheets.Item("1_8")
$old.Delete() | Out-Null

# Add the new "6.Raw" column header, matching the bold/bordered header style
# already used by B1:F1.
$ws.Range("G1").Value = "6.Raw"
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)

# Row 2 becomes the single "Test" capture row with fresh DMM readings.
$ws.Range("A2").Value = "Test"
$ws.Range("B2").Value = -7.8342
$ws.Range("C2").Value = -7.8254
$ws.Range("D2").Value = -7.8416
$ws.Range("E2").Value = 0.005
$ws.Range("F2").Value = 10
$ws.Range("G2").Value = "-7.836767,-7.83968,-7.841563,-7.833445,-7.825357,-7.831563,-7.82865,-7.837141,-7.837677,-7.829997"

# Drop the old "Deep_Sleep" row (row 3) that's no longer needed.
$ws.Rows.Item(3).Delete() | Out-Null

# Rename the surviving sheet.
$ws.Name = "test"
